$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5619495511054993
$ws.Range("B1").Value = 1.479623913764954
$ws.Range("C1").Value = 4.54067325592041
$ws.Range("D1").Value = 1.422715067863464
$ws.Range("E1").Value = 0.8134395480155945
